$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: 2021年 -- copy formatting from the row above (A6) for the year-label cell
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = "2021年"

$ws.Range("B7").Value = 98.8
$ws.Range("C7").Value = 101.1
$ws.Range("D7").Value = 99.90000000000001
$ws.Range("E7").Value = 102.8
$ws.Range("F7").Value = 100.8
$ws.Range("G7").Value = 99.90000000000001
$ws.Range("H7").Value = 100.4

# Row 8: 2022年 -- only H populated; B:G left blank (empty strings), like the diff
$ws.Range("A6").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "2022年"

$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = 101.2
